$d = $word.ActiveDocument
$nl = [char]11

function Set-ParaText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $r.Text = $newText
}

# --- Title ---
Set-ParaText 1 "Unveiling the Enigma of Chemistry: A Journey into the Realm of Matter and Transformations"

# --- Author name: "Emily Taylor" -> "Dr. Avery Donovan" ---
Set-ParaText 2 "Dr. Avery Donovan"

# --- Email: "etaylor@biodiscovery.org" -> "chemistry.dr.donovan@highschool.academy" ---
Set-ParaText 3 "chemistry.dr.donovan@highschool.academy"

# --- Main body paragraph (paragraph 5, 1-indexed) ---
$body = "In the vast tapestry of sciences, chemistry stands as a beacon of discovery, illuminating the composition of matter and the intricate dance of transformations it undergoes." `
    + " As we embark on this captivating journey into the realm of chemistry, we will unravel the secrets of the atom, witness the symphony of chemical reactions, and explore the profound implications of chemistry in medicine, industry, and everyday life." `
    + $nl + $nl `
    + "Chemistry unveils the enigmatic world of particles, revealing the intricate dance of atoms, ions, and molecules." `
    + " We will delve into the depths of the periodic table, deciphering the patterns and properties that govern the elements." `
    + " From the fiery brilliance of lithium to the noble elegance of helium, each element holds a unique story, waiting to be explored." `
    + $nl + $nl `
    + "The macroscopic world is a stage on which chemistry plays a transformative role." `
    + " Chemical reactions, like choreographed ballets, orchestrate spectacular displays of color, energy, and matter." `
    + " From the explosive combustion of fuels to the gentle rusting of iron, chemical reactions shape our world in countless ways." `
    + $nl + $nl `
    + "Beyond the theoretical realm, chemistry finds boundless applications in medicine, industry, and everyday life." `
    + " From life-saving drugs to durable materials, chemistry touches every aspect of our existence." `
    + " We will uncover the myriad ways in which chemistry contributes to our health, well-being, and technological advancements."
Set-ParaText 5 $body

# --- Summary heading stays "Summary" (paragraph 6) - unchanged ---

# --- Summary body paragraph (paragraph 7, 1-indexed) ---
$summary = "In this essay, we embarked on a captivating expedition into the world of chemistry, delving into the microscopic realm of particles and exploring the enchanting choreography of chemical reactions." `
    + " We discovered the profound applications of chemistry in medicine, industry, and everyday life." `
    + " Chemistry, with its ability to unravel the secrets of matter and orchestrate transformations, stands as a pillar of scientific knowledge, enriching our understanding of the universe and empowering us to create a better future."
Set-ParaText 7 $summary

# --- Add a new empty paragraph at the very end of the document body ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# --- Fix font name everywhere: "TimesNewToman" -> "Times New Roman" ---
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.End -gt $r.Start) {
        $r2 = $d.Range($r.Start, $r.End - 1)
        if ($r2.End -gt $r2.Start) {
            $r2.Font.Name = "Times New Roman"
        }
    }
}

Write-Output "done"
